$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Models")
$ws.Range("A1").Value = "Test"
